$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin table: updated prices / 1h volume deltas for almost
# every row, plus Uniswap (row 19) and Chainlink (row 20) swapping places.
#
# Every value here is plain text in the workbook (even "price" cells like
# "1.00" or "11.37"), so a leading apostrophe forces Excel to store each
# assignment as literal text instead of auto-converting number-looking
# strings into real numbers. The Style reset afterwards clears the
# "quote prefix" formatting flag Excel attaches in that case, so the cell
# ends up with the same (default) style as before the edit.
$ws.Range('D2').Value = "'" + "70.153.54"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'" + "  -3.04%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'" + "2.521.18"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'" + "  -4.29%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'" + "  +0.04%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'" + "578.22"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'" + "  -1.28%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'" + "168.38"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'" + "  -3.88%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'" + "  +0.09%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'" + "  -0.12%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'" + "2.521.12"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'" + "  -4.22%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'" + "  -5.84%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'" + "  -1.63%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'" + "0.349"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'" + "  -2.69%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'" + "4.91"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'" + "  -0.35%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'" + "2.983.98"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'" + "  -4.20%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'" + "70.040.01"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'" + "  -2.97%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'" + "  -5.33%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'" + "25.17"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'" + "  -2.19%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'" + "2.531.51"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'" + "  -3.74%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').Value = "'" + "Chainlink"
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').Value = "'" + "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').Value = "'" + "11.37"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'" + "  -5.55%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').Value = "'" + "Uniswap"
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').Value = "'" + "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').Value = "'" + "7.79"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'" + "  -0.89%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'" + "351.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'" + "  -6.33%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'" + "  -3.84%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'" + "1.96"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'" + "  -3.67%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'" + "  +0.19%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'" + "69.31"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'" + "  -3.06%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'" + "4.01"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'" + "  -5.05%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'" + "9.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'" + "  -4.36%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'" + "  -4.36%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'" + "0.999"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'" + "  +0.51%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'" + "  -3.93%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'" + "7.92"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'" + "  -0.34%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'" + "  -2.55%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'" + "466.47"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'" + "  -4.99%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D35').Value = "'" + "1.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'" + "  +0.12%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'" + "  +3.71%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'" + "153.24"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'" + "  -5.07%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'" + "  +0.72%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'" + "18.50"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'" + "  -3.40%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'" + "  +0.06%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'" + "4.79"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'" + "  -1.77%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'" + "  -1.32%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'" + "  -6.96%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'" + "  -13.66%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'" + "2.32"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'" + "  -9.38%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'" + "38.19"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'" + "  -2.16%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'" + "143.54"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'" + "  -4.47%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'" + "0.533"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'" + "  -1.58%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'" + "  -3.13%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'" + "  -4.23%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'" + "  -1.19%  "
$ws.Range('E51').Style = 'Normal'
